# Hortaliza, Macroferia Regional de Talca - Zapallo
# A new weekly price record was added to the top of the Camote data block
# (row 396), shifting every subsequent record down by one row and growing
# the sheet from 458 to 459 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 396:458 down to 397:459, leaving a blank row 396 to populate.
$ws.Rows("396:396").Insert()

# Fill the newly inserted row with the new record's data.
$ws.Cells.Item(396, 1).Value = 5
$ws.Cells.Item(396, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(396, 3).Value = "Maule"
$ws.Cells.Item(396, 4).Value = 45127
$ws.Cells.Item(396, 5).Value = 7
$ws.Cells.Item(396, 6).Value = 100112045
$ws.Cells.Item(396, 7).Value = "Zapallo"
$ws.Cells.Item(396, 8).Value = "Camote"
$ws.Cells.Item(396, 9).Value = "1a (guarda)"
$ws.Cells.Item(396, 10).Value = 800
$ws.Cells.Item(396, 11).Value = 400
$ws.Cells.Item(396, 12).Value = 400
$ws.Cells.Item(396, 13).Value = 400
$ws.Cells.Item(396, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(396, 15).Value = "Región del Maule"
$ws.Cells.Item(396, 16).Value = 400
$ws.Cells.Item(396, 17).Value = 1
$ws.Cells.Item(396, 18).Value = "Hortaliza"
